# Edit: Included ZW General Insurance and Lending business lines.
$wb = $excel.ActiveWorkbook

# --- 1. Unhide the "valid_product_names" sheet ---
$validSheet = $wb.Worksheets.Item("valid_product_names")
$validSheet.Visible = -1

# --- 2. Populate the Zimbabwe product list (column B, rows 2-68) ---
$zimbabweProducts = @(
    "Flexi Funeral Plan"
    "Agribusiness Loan"
    "Working Capital Loans- Group Clients"
    "Order Finance Loans"
    "Working Capital Loans- Individual clients (Unsecured loans)"
    "Salary Based Loan"
    "Structured Finance Loan"
    "SOFL"
    "SFN"
    "Savings Plan"
    "Life Plan"
    "Funeral Plan"
    "Term Plan"
    "Equity - General"
    "Fixed Interest - Money Market"
    "Real Estate - General"
    "HOME PLAN"
    "Private Motor Car"
    "Personal Package"
    "PTA Yellow Card"
    "Motor Combined"
    "Motor Act"
    "Motor Fleet"
    "Business Package"
    "Employers Personal Accident"
    "Assets All Risks"
    "Employers/Residual Liability"
    "Fronting Product"
    "Liabilities"
    "Machinery Breakdown"
    "Marine Cargo"
    "MOTOR PLAN"
    "Marine Hull"
    "Bonds - Court"
    "Personal Combined: Old Mutual"
    "CABS"
    "Fidelity Guarantee"
    "Houseowners"
    "Contractors"
    "Personal Accident"
    "Marine Open Policy"
    "Motor Fleet Eaton & Young"
    "Farmers Package"
    "Personal Combined: RMI"
    "Zimbabwe Caravan Association"
    "SME Business Package"
    "Kingsure Personal Package"
    "Fire"
    "Machinery Loss of Profits"
    "Travel Insurance"
    "Agrisure Personal Package"
    "Professional Indemnity"
    "Reinsurance For All Products"
    "Motor Traders External"
    "Livestock & Bloodstock"
    "MBCA Insure Motor Plan"
    "MBCA Insure Home Plan"
    "Mortgage Guarantee"
    "Forex Travel Insurance"
    "Aviation"
    "Stanchart Personal Pcakage"
    "Electronic Gadgets Insurance"
    "Allsure"
    "Emergency Rescue"
    "Medic-Sure"
    "Loss of Profits"
    "Living Future"
)

for ($i = 0; $i -lt $zimbabweProducts.Length; $i++) {
    $validSheet.Cells.Item($i + 2, 2).Value = $zimbabweProducts[$i]
}

# --- 3. Expand the "Zimbabwe" named range to cover the new rows ---
$wb.Names.Item("Zimbabwe").RefersTo = "=valid_product_names!`$B`$2:`$B`$68"

# --- 4. Update the "customer_prod" sheet sample/demo values ---
$prodSheet = $wb.Worksheets.Item("customer_prod")
$prodSheet.Range("A2").Value = "Equity - General"
$prodSheet.Range("C2").Value = "Zimbabwe"
$prodSheet.Range("A3").Value = "Funeral Plan"
$prodSheet.Range("A4").Value = "Salary Based Loan"
$prodSheet.Range("A5").Value = "Private Motor Car"

# --- 5. Restore selections to match the authored state ---
$validSheet.Activate()
$validSheet.Range("A1:A62,B1:B68").Select()
$validSheet.Range("B1").Activate()

$prodSheet.Activate()
$prodSheet.Range("C7").Select()
